# Refresh crypto price/volume snapshot (scheduled GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.167.48"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "2.092.13"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'228.51"
$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "'60.99"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "2.402.23"
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'22.28"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("D15").Value = "'5.48"
$ws.Range("D16").Value = "'0.776"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "2.106.64"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "38.113.93"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "'70.20"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'224.02"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").Value = "'170.04"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'18.97"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'1.36"
$ws.Range("E30").Value = "  +6.57%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +6.35%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  +5.65%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'18.09"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("D41").Value = "1.552.67"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "'100.13"
$ws.Range("E42").Value = "  +4.17%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "'0.0913"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "'7.46"
$ws.Range("E48").Value = "  +4.54%  "
$ws.Range("D49").Value = "'1.02"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "2.289.24"
$ws.Range("E51").Value = "  +2.90%  "
